$d = $word.ActiveDocument

# Helper: find the 1-based index of the first paragraph (at/after $startAt)
# whose full text (including its trailing paragraph mark) equals $text.
function Find-ParaIndex($doc, $text, $startAt) {
    $count = $doc.Paragraphs.Count
    for ($i = $startAt; $i -le $count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -eq $text) {
            return $i
        }
    }
    return -1
}

# Locate the entity headings (top-level list items) that scope each block,
# so same-named fields that repeat across entities aren't confused with one
# another (e.g. "ID музыканта" / "ID публикации" / "ID автора" all occur
# more than once in the document).
$hUser = Find-ParaIndex $d "Пользователь`r" 1
$hPost = Find-ParaIndex $d "Публикация`r" 1
$hMusician = Find-ParaIndex $d "Музыкант`r" 1
$hComment = Find-ParaIndex $d "Комментарий`r" 1
$hRating = Find-ParaIndex $d "Оценка`r" 1

# --- Сущность "Пользователь": underline the "ID пользователя" field ---
$fUserId = Find-ParaIndex $d "ID пользователя`r" $hUser
$d.Paragraphs.Item($fUserId).Range.Font.Underline = 1

# --- Сущность "Публикация": underline the "ID публикации" field ---
$fPostId = Find-ParaIndex $d "ID публикации`r" $hPost
$d.Paragraphs.Item($fPostId).Range.Font.Underline = 1

# --- Сущность "Музыкант": underline the "ID музыканта" field ---
$fMusicianId = Find-ParaIndex $d "ID музыканта`r" $hMusician
$d.Paragraphs.Item($fMusicianId).Range.Font.Underline = 1

# --- Сущность "Комментарий": underline the "ID комментария" field ---
$fCommentId = Find-ParaIndex $d "ID комментария`r" $hComment
$d.Paragraphs.Item($fCommentId).Range.Font.Underline = 1

# --- Сущность "Оценка": ---
# Resolve all three field paragraphs up front (before any text edits or
# deletions shift paragraph indices).
$fRatingId = Find-ParaIndex $d "ID оценки`r" $hRating
$fRatingPostId = Find-ParaIndex $d "ID публикации`r" $hRating
$fRatingAuthorId = Find-ParaIndex $d "ID автора`r" $hRating

# 1) "ID оценки" becomes "ID публикации" (underlined), keeping the "ID "
#    run (with its en-US language tag) intact and only retargeting the
#    word after it.
$p = $d.Paragraphs.Item($fRatingId)
$p.Range.Font.Underline = 1
$pStart = $p.Range.Start
$pEnd = $p.Range.End
$wordRange = $d.Range($pStart + 3, $pEnd - 1)
$wordRange.Text = "публикации"

# 2) "ID автора" stays as-is, just gets underlined.
$d.Paragraphs.Item($fRatingAuthorId).Range.Font.Underline = 1

# 3) The original "ID публикации" field is removed entirely (deleted last
#    so the indices resolved above stay valid).
$d.Paragraphs.Item($fRatingPostId).Range.Delete()

Write-Output "edits applied"
